$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.233.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7089"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08004"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.84%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3031"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08209"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.183"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.44%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.791.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.15%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7038"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.15%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.916.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.821"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007880"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.480"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "

# Row 24
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.949.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.913"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1443"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.917"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.424"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.478"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.362"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.021"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05192"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.38%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.671"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9307"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.138.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4255"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.880"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5340"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.769"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.168"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.956.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.08%  "
